$d = $word.ActiveDocument
$d.Content.Find.Execute("2023-04-23 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-24 Monday", 2) | Out-Null
$t = $d.Tables.Item(1)
$t.Cell(1,1).Range.Text = "18×36=648"
$t.Cell(1,2).Range.Text = "90×91=8190"
$t.Cell(1,3).Range.Text = "76×18=1368"
$t.Cell(1,4).Range.Text = "72×82=5904"
$t.Cell(1,5).Range.Text = "46×76=3496"
$t.Cell(2,1).Range.Text = "19×97=1843"
$t.Cell(2,2).Range.Text = "62×94=5828"
$t.Cell(2,3).Range.Text = "86×15=1290"
$t.Cell(2,4).Range.Text = "51×96=4896"
$t.Cell(2,5).Range.Text = "44×70=3080"
$t.Cell(3,1).Range.Text = "61×62=3782"
$t.Cell(3,2).Range.Text = "76×56=4256"
$t.Cell(3,3).Range.Text = "50×92=4600"
$t.Cell(3,4).Range.Text = "12×58=696"
$t.Cell(3,5).Range.Text = "34×51=1734"
$t.Cell(4,1).Range.Text = "88×72=6336"
$t.Cell(4,2).Range.Text = "98×20=1960"
$t.Cell(4,3).Range.Text = "23×10=230"
$t.Cell(4,4).Range.Text = "71×41=2911"
$t.Cell(4,5).Range.Text = "97×53=5141"
$t.Cell(5,1).Range.Text = "87×81=7047"
$t.Cell(5,2).Range.Text = "73×96=7008"
$t.Cell(5,3).Range.Text = "10×76=760"
$t.Cell(5,4).Range.Text = "17×11=187"
$t.Cell(5,5).Range.Text = "77×98=7546"
$t.Cell(6,1).Range.Text = "32×76=2432"
$t.Cell(6,2).Range.Text = "37×54=1998"
$t.Cell(6,3).Range.Text = "41×96=3936"
$t.Cell(6,4).Range.Text = "93×34=3162"
$t.Cell(6,5).Range.Text = "93×62=5766"
$t.Cell(7,1).Range.Text = "45×99=4455"
$t.Cell(7,2).Range.Text = "80×29=2320"
$t.Cell(7,3).Range.Text = "71×80=5680"
$t.Cell(7,4).Range.Text = "67×43=2881"
$t.Cell(7,5).Range.Text = "24×29=696"
$t.Cell(8,1).Range.Text = "82×36=2952"
$t.Cell(8,2).Range.Text = "78×33=2574"
$t.Cell(8,3).Range.Text = "18×75=1350"
$t.Cell(8,4).Range.Text = "61×39=2379"
$t.Cell(8,5).Range.Text = "98×94=9212"
$t.Cell(9,1).Range.Text = "58×11=638"
$t.Cell(9,2).Range.Text = "89×39=3471"
$t.Cell(9,3).Range.Text = "95×48=4560"
$t.Cell(9,4).Range.Text = "18×92=1656"
$t.Cell(9,5).Range.Text = "55×27=1485"
$t.Cell(10,1).Range.Text = "77×41=3157"
$t.Cell(10,2).Range.Text = "76×25=1900"
$t.Cell(10,3).Range.Text = "62×47=2914"
$t.Cell(10,4).Range.Text = "21×63=1323"
$t.Cell(10,5).Range.Text = "71×25=1775"
$t.Cell(11,1).Range.Text = "38×75=2850"
$t.Cell(11,2).Range.Text = "21×95=1995"
$t.Cell(11,3).Range.Text = "89×48=4272"
$t.Cell(11,4).Range.Text = "38×13=494"
$t.Cell(11,5).Range.Text = "13×58=754"
$t.Cell(12,1).Range.Text = "39×94=3666"
$t.Cell(12,2).Range.Text = "19×33=627"
$t.Cell(12,3).Range.Text = "18×67=1206"
$t.Cell(12,4).Range.Text = "16×86=1376"
$t.Cell(12,5).Range.Text = "38×35=1330"
$t.Cell(13,1).Range.Text = "47×72=3384"
$t.Cell(13,2).Range.Text = "13×58=754"
$t.Cell(13,3).Range.Text = "38×48=1824"
$t.Cell(13,4).Range.Text = "24×49=1176"
$t.Cell(13,5).Range.Text = "62×81=5022"
$t.Cell(14,1).Range.Text = "96×47=4512"
$t.Cell(14,2).Range.Text = "98×40=3920"
$t.Cell(14,3).Range.Text = "72×73=5256"
$t.Cell(14,4).Range.Text = "78×24=1872"
$t.Cell(14,5).Range.Text = "37×21=777"
$t.Cell(15,1).Range.Text = "35×13=455"
$t.Cell(15,2).Range.Text = "47×34=1598"
$t.Cell(15,3).Range.Text = "55×73=4015"
$t.Cell(15,4).Range.Text = "87×36=3132"
$t.Cell(15,5).Range.Text = "69×75=5175"
$t.Cell(16,1).Range.Text = "27×93=2511"
$t.Cell(16,2).Range.Text = "55×15=825"
$t.Cell(16,3).Range.Text = "55×67=3685"
$t.Cell(16,4).Range.Text = "72×39=2808"
$t.Cell(16,5).Range.Text = "15×45=675"
$t.Cell(17,1).Range.Text = "57×89=5073"
$t.Cell(17,2).Range.Text = "24×65=1560"
$t.Cell(17,3).Range.Text = "97×16=1552"
$t.Cell(17,4).Range.Text = "74×74=5476"
$t.Cell(17,5).Range.Text = "68×55=3740"
$t.Cell(18,1).Range.Text = "73×22=1606"
$t.Cell(18,2).Range.Text = "55×28=1540"
$t.Cell(18,3).Range.Text = "86×80=6880"
$t.Cell(18,4).Range.Text = "20×13=260"
$t.Cell(18,5).Range.Text = "76×94=7144"
$t.Cell(19,1).Range.Text = "87×91=7917"
$t.Cell(19,2).Range.Text = "99×83=8217"
$t.Cell(19,3).Range.Text = "27×53=1431"
$t.Cell(19,4).Range.Text = "86×58=4988"
$t.Cell(19,5).Range.Text = "56×81=4536"
$t.Cell(20,1).Range.Text = "42×87=3654"
$t.Cell(20,2).Range.Text = "42×66=2772"
$t.Cell(20,3).Range.Text = "62×48=2976"
$t.Cell(20,4).Range.Text = "32×49=1568"
$t.Cell(20,5).Range.Text = "44×46=2024"
